$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D16").Value = "2016-03-04 11:15:32"
$wsZhCn.Range("G16").Value = "2016-03-04 11:16:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-03-04 11:15:57"
$wsDeDe.Range("G5").Value = "2016-03-04 11:17:18"
